$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - "Organising/enter data for database and images into res folder for the glossary."
# Status moves from "In Progress" to "Done"; Tuesday hours reduced to 0, add Wed (L10) = 0
$ws.Range("D10").Value = "Done"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

# Row 45 - "Unit tests and proof of them passing"
# Assign to Kari, mark as Done, and fill in remaining day columns with 0
$ws.Range("C45").Value = "Kari"
$ws.Range("D45").Value = "Done"
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0

# Row 48 - "Burndown chart"
# Remove the stray value in column C (Who) and shift the day value into K48 instead
$ws.Range("C48").ClearContents()
$ws.Range("K48").Value = 2

# Row 49 - "Refactioring"
# Remove the stray value in column C (Who) and shift the day value into K49 instead
$ws.Range("C49").ClearContents()
$ws.Range("K49").Value = 6

# Update the view so it reflects scrolling to row 26 with N49 selected
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("N49").Select()
